$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump Version and Date for the 2.0.0 terminology IG publish
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.0.1"

# The new date string ("2025-09-22") looks like a real date to Excel's literal
# parser, so assigning it directly would silently convert the cell to a date
# serial number. Force it in as literal text (leading apostrophe) and then
# restore the plain body-row style (copied from a sibling text cell) so the
# cell keeps looking/behaving like every other "Value" cell on this sheet.
$meta.Range("B8").Value = "'2025-09-22"
$meta.Range("B2").Copy()
$meta.Range("B8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Split the old "Concepts" sheet into two sheets:
#      - "Properties" (reuses the old sheet/rId/sheetId, just renamed) holding
#        the CodeSystem property definitions (status / effectiveDate)
#      - a brand-new "Concepts" sheet holding the concept list that used to
#        live on the sheet we are about to rename.
# ---------------------------------------------------------------------------
$concepts = $wb.Worksheets.Item("Concepts")

# Create the new sheet right after the current "Concepts" sheet and copy the
# existing concept rows (values + styles) into it before we repurpose the
# original sheet as "Properties".
$newConcepts = $wb.Worksheets.Add($null, $concepts)
$concepts.Range("A1:D3").Copy($newConcepts.Range("A1"))

# Rename the original sheet to "Properties" and the freshly added one to
# "Concepts", so the tab order becomes Metadata, Properties, Concepts.
$concepts.Name = "Properties"
$newConcepts.Name = "Concepts"

# ---------------------------------------------------------------------------
# 3. Populate the "Properties" sheet with the CodeSystem property table.
# ---------------------------------------------------------------------------
$props = $wb.Worksheets.Item("Properties")

$props.Range("A1").Value = "Code"
$props.Range("B1").Value = "Uri"
$props.Range("C1").Value = "Description"
$props.Range("D1").Value = "Type"

$props.Range("A2").Value = "status"
$props.Range("B2").Value = "http://hl7.org/fhir/concept-properties#status"
$props.Range("C2").Value = "A property that indicates the status of the concept. One of active, experimental, deprecated, or retired."
$props.Range("D2").Value = "code"

$props.Range("A3").Value = "effectiveDate"
$props.Range("B3").Value = "http://hl7.org/fhir/concept-properties#effectiveDate"
$props.Range("C3").Value = "The date at which the concept status was last changed."
$props.Range("D3").Value = "dateTime"
